# 17.6.1 Fixed Internet broadband subscriptions by speed*
# - capitalise the footnote ("*according..." -> "*According...")
# - add a new "2023" column (O) with its data
# - bump the header / footnote row heights to fit
# - shrink the footnote row's font so it still fits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix capitalisation of the English footnote in C7
# ---------------------------------------------------------------------
$ws.Range("C7").Value = "*According to the Service for the Regulation and Supervision of the Communications Sector under the Ministry of Digital Development of the Kyrgyz Republic"

# ---------------------------------------------------------------------
# 2. Add the new "2023" column (O) of data, copying the neighbouring
#    (N) column's number formatting / borders for each row.
# ---------------------------------------------------------------------

# Row 2 - thin separator row below the title, just needs the same
# bottom border formatting as the rest of the row.
$ws.Range("N2").Copy() | Out-Null
$ws.Range("O2").PasteSpecial(-4122) | Out-Null

# Row 3 - year headers
$ws.Range("O3").Value = 2023
$ws.Range("N3").Copy() | Out-Null
$ws.Range("O3").PasteSpecial(-4122) | Out-Null

# Row 4 - "256 kbit/s to < 2 Mbit/s" data
$ws.Range("O4").Value = 5571
$ws.Range("N4").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null

# Row 5 - "2 Mbit/s to < 10 Mbit/s" data
$ws.Range("O5").Value = 74710
$ws.Range("N5").Copy() | Out-Null
$ws.Range("O5").PasteSpecial(-4122) | Out-Null

# Row 6 - ">= 10 Mbit/s" data (bottom thick border row)
$ws.Range("O6").Value = 375715
$ws.Range("N6").Copy() | Out-Null
$ws.Range("O6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Row-height tweaks to accommodate the wider table
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 46.5

# ---------------------------------------------------------------------
# 4. Shrink the footnote row (A7:C7) font so the row still reads well
#    now that there is an extra column.
# ---------------------------------------------------------------------
$ws.Range("A7:C7").Font.Size = 8
